# Apply the changes described by the diff:
#  1. Insert a new "Author" paragraph ("Dr. Lennart Wittkuhn") at the very
#     start of the document body.
#  2. Set the first table's layout to "fixed" (adds <w:tblLayout w:type="fixed"/>).
#  3. Add a new paragraph style "AbstractTitle" ("Abstract Title").
#  4. Change the existing "Abstract" style's space-before from 300 -> 100 (twips).
#  5. Add a new paragraph style "FootnoteBlockText" ("Footnote Block Text").

$d = $word.ActiveDocument

# --- 1. New Author paragraph at the top of the document -------------------
$introRange = $d.Range(0, 0)
$introRange.InsertParagraphBefore()

$authorPara = $d.Paragraphs.Item(1)
$authorPara.Style = "Author"

$authorRange = $authorPara.Range
$authorRange.End = $authorRange.End - 1
$authorRange.Text = "Dr. Lennart Wittkuhn"

# --- 2. Fixed table layout on the first (only) table -----------------------
$table = $d.Tables.Item(1)
$table.AutoFitBehavior(0)

# --- 3. New "AbstractTitle" style -------------------------------------------
$abstractTitle = $d.Styles.Add("AbstractTitle", 1)
$abstractTitle.NameLocal = "Abstract Title"
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true

$atPf = $abstractTitle.ParagraphFormat
$atPf.KeepWithNext = $true
$atPf.KeepTogether = $true
$atPf.Alignment = 1
$atPf.SpaceBefore = 15
$atPf.SpaceAfter = 0

$atFont = $abstractTitle.Font
$atFont.Size = 10
$atFont.SizeBi = 10
$atFont.Bold = $true
$atFont.Color = 9067060

# --- 4. Abstract style spacing tweak ---------------------------------------
$abstract = $d.Styles.Item("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# --- 5. New "FootnoteBlockText" style ---------------------------------------
$footnoteBlock = $d.Styles.Add("FootnoteBlockText", 1)
$footnoteBlock.NameLocal = "Footnote Block Text"
$footnoteBlock.BaseStyle = "Footnote Text"
$footnoteBlock.NextParagraphStyle = "Footnote Text"
$footnoteBlock.Priority = 9
$footnoteBlock.UnhideWhenUsed = $true
$footnoteBlock.QuickStyle = $true

$fbPf = $footnoteBlock.ParagraphFormat
$fbPf.SpaceBefore = 5
$fbPf.SpaceAfter = 5
$fbPf.LeftIndent = 24
$fbPf.RightIndent = 24
$fbPf.FirstLineIndent = 0

Write-Host "Edit applied."
